$wb = $excel.ActiveWorkbook

# --- Regular cell value updates (data refresh for 2023-04-17) ---
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 44
$ws.Range("E3").Value = 42
$ws.Range("G3").Value = 32
$ws.Range("I3").Value = 50
$ws.Range("D6").Value = 128
$ws.Range("H6").Value = 115
$ws.Range("J6").Value = 118
$ws.Range("D7").Value = 192
$ws.Range("E7").Value = 195
$ws.Range("G7").Value = 228
$ws.Range("H7").Value = 176
$ws.Range("I7").Value = 236
$ws.Range("J7").Value = 222

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 7
$ws.Range("G6").Value = 6
$ws.Range("I6").Value = 15

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("H5").Value = 3
$ws.Range("H6").Value = 5

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J17").Value = 7
$ws.Range("J29").Value = 1
$ws.Range("G34").Value = 6
$ws.Range("I34").Value = 15
$ws.Range("H48").Value = 5
$ws.Range("H51").Value = 22
$ws.Range("I51").Value = 39
$ws.Range("H60").Value = 3
$ws.Range("D67").Value = 4
$ws.Range("E72").Value = 7
$ws.Range("H77").Value = 4
$ws.Range("D93").Value = 192
$ws.Range("E93").Value = 195
$ws.Range("G93").Value = 228
$ws.Range("H93").Value = 176
$ws.Range("I93").Value = 236
$ws.Range("J93").Value = 222

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I3").Value = 9
$ws.Range("H6").Value = 17
$ws.Range("H7").Value = 22
$ws.Range("I7").Value = 39

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("D3").Value = 4
$ws.Range("D4").Value = 4

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("E3").Value = 4
$ws.Range("E6").Value = 7

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 3

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 3
$ws.Range("J5").Value = 7

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("B4").Value = 3
$ws.Range("B5").Value = 4
# --- Galewood: add new 2023 column (D) ---
$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("C1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = 2023
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth
